$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-CellText $ws "D2" "43.401.82"
Set-CellText $ws "E2" "  +1.29%  "
Set-CellText $ws "D3" "2.235.12"
Set-CellText $ws "E3" "  +0.13%  "
Set-CellText $ws "E4" "  -0.25%  "
Set-CellText $ws "D5" "317.73"
Set-CellText $ws "E5" "  +1.37%  "
Set-CellText $ws "D6" "99.41"
Set-CellText $ws "E6" "  -0.33%  "
Set-CellText $ws "E7" "  +2.12%  "
Set-CellText $ws "E8" "  -0.12%  "
Set-CellText $ws "D9" "0.563"
Set-CellText $ws "E9" "  +0.41%  "
Set-CellText $ws "D10" "37.26"
Set-CellText $ws "E10" "  -0.27%  "
Set-CellText $ws "D11" "0.0830"
Set-CellText $ws "E11" "  -0.97%  "
Set-CellText $ws "D12" "7.70"
Set-CellText $ws "E12" "  +1.42%  "
Set-CellText $ws "E13" "  +1.80%  "
Set-CellText $ws "B14" "WrappedliquidstakedEther2.0"
Set-CellText $ws "C14" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-CellText $ws "D14" "2.575.84"
Set-CellText $ws "E14" "  +0.11%  "
Set-CellText $ws "B15" "Polygon"
Set-CellText $ws "C15" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-CellText $ws "D15" "0.866"
Set-CellText $ws "E15" "  -1.30%  "
Set-CellText $ws "D16" "14.39"
Set-CellText $ws "E16" "  +4.01%  "
Set-CellText $ws "D17" "2.219.90"
Set-CellText $ws "E17" "  -1.03%  "
Set-CellText $ws "D18" "43.368.97"
Set-CellText $ws "E18" "  +1.72%  "
Set-CellText $ws "D19" "14.19"
Set-CellText $ws "E19" "  -1.08%  "
Set-CellText $ws "D20" "6.63"
Set-CellText $ws "E20" "  +0.13%  "
Set-CellText $ws "D21" "0.0₃0971"
Set-CellText $ws "E21" "  +2.38%  "
Set-CellText $ws "D22" "3.22"
Set-CellText $ws "E22" "  -1.57%  "
Set-CellText $ws "D23" "65.27"
Set-CellText $ws "E23" "  +0.47%  "
Set-CellText $ws "D24" "236.25"
Set-CellText $ws "E24" "  +0.68%  "
Set-CellText $ws "E25" "  +2.43%  "
Set-CellText $ws "E26" "  +0.18%  "
Set-CellText $ws "E27" "  +2.83%  "
Set-CellText $ws "D28" "10.10"
Set-CellText $ws "E28" "  -1.36%  "
Set-CellText $ws "E29" "  +2.30%  "
Set-CellText $ws "D30" "6.40"
Set-CellText $ws "E30" "  -3.24%  "
Set-CellText $ws "D31" "36.76"
Set-CellText $ws "E31" "  +9.50%  "
Set-CellText $ws "B32" "EthereumClassic"
Set-CellText $ws "C32" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-CellText $ws "D32" "20.31"
Set-CellText $ws "E32" "  -1.30%  "
Set-CellText $ws "B33" "Hedera"
Set-CellText $ws "C33" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-CellText $ws "D33" "0.0875"
Set-CellText $ws "E33" "  -1.57%  "
Set-CellText $ws "D34" "157.75"
Set-CellText $ws "E34" "  -2.78%  "
Set-CellText $ws "D35" "2.71"
Set-CellText $ws "E35" "  -2.07%  "
Set-CellText $ws "D36" "3.24"
Set-CellText $ws "E36" "  +5.95%  "
Set-CellText $ws "E37" "  +0.02%  "
Set-CellText $ws "E38" "  +0.81%  "
Set-CellText $ws "D39" "4.42"
Set-CellText $ws "E39" "  +0.40%  "
Set-CellText $ws "D40" "0.104"
Set-CellText $ws "E40" "  -1.17%  "
Set-CellText $ws "D41" "3.70"
Set-CellText $ws "E41" "  +1.85%  "
Set-CellText $ws "D42" "0.0322"
Set-CellText $ws "E42" "  -0.08%  "
Set-CellText $ws "D43" "14.48"
Set-CellText $ws "E43" "  +19.54%  "
Set-CellText $ws "E44" "  +0.01%  "
Set-CellText $ws "D45" "1.807.97"
Set-CellText $ws "E45" "  -0.51%  "
Set-CellText $ws "D46" "0.203"
Set-CellText $ws "E46" "  -1.70%  "
Set-CellText $ws "D47" "84.53"
Set-CellText $ws "E47" "  -6.26%  "
Set-CellText $ws "E48" "  -0.98%  "
Set-CellText $ws "D49" "8.84"
Set-CellText $ws "E49" "  +1.77%  "
Set-CellText $ws "D50" "74.38"
Set-CellText $ws "E50" "  -2.75%  "
Set-CellText $ws "D51" "58.72"
Set-CellText $ws "E51" "  -3.83%  "
